$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# December ("desember") row 13 was only partially filled in (label only).
# Add the BPJS/Umum counts and the Total formula, matching the pattern
# already used by the rows above (Total = BPJS + Umum).
$ws.Range("B13").Value = 67
$ws.Range("C13").Value = 24
$ws.Range("D13").Formula = "=SUM(B13:C13)"

# The saved view's active cell/selection moved to D18.
$ws.Range("D18").Select()
